$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $a, $b, $c, $d) {
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $a
    $ws.Range("A$r").Style = "Normal"
    $ws.Range("B$r").Value = $b
    $ws.Range("C$r").Value = $c
    $ws.Range("D$r").Value = $d
}

Set-Row 623 "2026/01/12" "月" 22 201
Set-Row 624 "2026/01/13" "火" 1 201
Set-Row 625 "2026/12/29" "火" 13 201
Set-Row 626 "2026/12/29" "火" 16 201
Set-Row 627 "2026/12/29" "火" 19 201
Set-Row 628 "2026/12/29" "火" 23 201
Set-Row 629 "2026/12/30" "水" 2 201
Set-Row 630 "2026/12/30" "水" 5 201
Set-Row 631 "2026/12/30" "水" 8 201
Set-Row 632 "2026/12/30" "水" 13 201
Set-Row 633 "2026/12/30" "水" 16 201
Set-Row 634 "2026/12/30" "水" 22 201
Set-Row 635 "2026/12/31" "木" 2 201
Set-Row 636 "2026/12/31" "木" 6 201
Set-Row 637 "2026/12/31" "木" 10 201
Set-Row 638 "2026/12/31" "木" 12 201
Set-Row 639 "2026/12/31" "木" 14 201
Set-Row 640 "2026/12/31" "木" 22 201
Set-Row 641 "2027/01/01" "金" 2 201
Set-Row 642 "2027/01/01" "金" 5 201
Set-Row 643 "2027/01/01" "金" 13 201
Set-Row 644 "2027/01/01" "金" 16 201
Set-Row 645 "2027/01/01" "金" 19 201
Set-Row 646 "2027/01/02" "土" 1 201
Set-Row 647 "2027/01/02" "土" 5 201
Set-Row 648 "2027/01/02" "土" 8 201
Set-Row 649 "2027/01/02" "土" 13 201
Set-Row 650 "2027/01/02" "土" 16 201
Set-Row 651 "2027/01/02" "土" 19 201
Set-Row 652 "2027/01/02" "土" 22 201
Set-Row 653 "2027/01/03" "日" 1 201
Set-Row 654 "2027/01/03" "日" 4 201
Set-Row 655 "2027/01/03" "日" 7 201
Set-Row 656 "2027/01/03" "日" 13 201
Set-Row 657 "2027/01/03" "日" 16 201
Set-Row 658 "2027/01/03" "日" 19 201
Set-Row 659 "2027/01/03" "日" 22 201
Set-Row 660 "2027/01/04" "月" 2 201
Set-Row 661 "2027/01/04" "月" 4 201
Set-Row 662 "2027/01/04" "月" 7 201
Set-Row 663 "2027/01/04" "月" 13 201
Set-Row 664 "2027/01/04" "月" 22 201
Set-Row 665 "2027/01/05" "火" 2 201
Set-Row 666 "2027/01/05" "火" 7 201
